# Generate Report for handoff
# Update "Latest Handoff Datetime" (column D) for rows whose status is
# "Handback transform failed" or "Ready for handoff" on both the zh-cn
# and de-de worksheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D4").Value  = "2016-02-16 15:57:05"
$zhcn.Range("D6").Value  = "2016-02-16 15:57:05"
$zhcn.Range("D7").Value  = "2016-02-16 15:57:05"
$zhcn.Range("D8").Value  = "2016-02-16 15:57:05"
$zhcn.Range("D9").Value  = "2016-02-16 15:57:05"
$zhcn.Range("D10").Value = "2016-02-16 15:57:05"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D4").Value  = "2016-02-16 15:57:22"
$dede.Range("D6").Value  = "2016-02-16 15:57:22"
$dede.Range("D7").Value  = "2016-02-16 15:57:22"
$dede.Range("D8").Value  = "2016-02-16 15:57:22"
$dede.Range("D9").Value  = "2016-02-16 15:57:22"
$dede.Range("D10").Value = "2016-02-16 15:57:22"
